# Generate Report for Handback
# Applies the localization-status "handback" update:
#  - Status text "Ready for handoff" -> "Handed back: in sync with en-US"
#    (Overview!E/F and the Status column on each language sheet)
#  - Populates "Latest Target File" / "Latest Handback File" /
#    "Latest Handback DateTime" columns on the zh-cn and de-de sheets
#    with the handed-back .md source hyperlink, the generated .xlf file,
#    and the handback timestamp.
#  - Widens a handful of columns to fit the newly-populated content.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$mdBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/06a15762ae7f88aec57d2f2e3cc2bdb7f09ed630/e2e/"

$file1 = "859e9e9f-dcc1-4d34-9199-9ef345bb5a9b"
$file2 = "c70a462c-d133-44c2-9e86-4df3f6cc1309"

# ---------------------------------------------------------------------
# Overview sheet: update the per-language status cells (E/F, rows 2-3)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Status column
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

# Row 2 (859e9e9f...)
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), ($mdBase + $file1 + ".md"), "", "", ($file1 + ".md"))
$wsZh.Range("J2").Value = ($file1 + ".7e7dfbe0f84d2017e4fc0d48a5f5c419ef44971a.zh-cn.xlf")
$wsZh.Range("K2").Value = "2016-08-16 11:02:31"

# Row 3 (c70a462c...)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), ($mdBase + $file2 + ".md"), "", "", ($file2 + ".md"))
$wsZh.Range("J3").Value = ($file2 + ".a669482bb75e2c896168956a5a11bb1c33b09a5e.zh-cn.xlf")
$wsZh.Range("K3").Value = "2016-08-16 11:02:31"

$wsZh.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZh.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZh.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Status column
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# Row 2 (859e9e9f...)
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), ($mdBase + $file1 + ".md"), "", "", ($file1 + ".md"))
$wsDe.Range("J2").Value = ($file1 + ".7e7dfbe0f84d2017e4fc0d48a5f5c419ef44971a.de-de.xlf")
$wsDe.Range("K2").Value = "2016-08-16 11:02:38"

# Row 3 (c70a462c...)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), ($mdBase + $file2 + ".md"), "", "", ($file2 + ".md"))
$wsDe.Range("J3").Value = ($file2 + ".a669482bb75e2c896168956a5a11bb1c33b09a5e.de-de.xlf")
$wsDe.Range("K3").Value = "2016-08-16 11:02:38"

$wsDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDe.Columns.Item(10).ColumnWidth = 39.166666666666664
